# Additional seeds for Displaced Jets test
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 7) down onto the
# four new rows so they pick up the same cell style without bloating the
# workbook's style table.
$ws.Range("A7:K7").Copy()
$ws.Range("A8:K11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$names = @("L1_SingleJetDISP15", "L1_SingleJetDISP30", "L1_SingleJetDISP45", "L1_SingleJetDISP90")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 1).Value = 6 + $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = 0
    for ($col = 4; $col -le 11; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}

$ws.Range("B13").Select() | Out-Null
